$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '51.905.31'
$ws.Range("E2").Value = '  +0.32%  '
$ws.Range("D3").Value = '2.815.53'
$ws.Range("E3").Value = '  +1.43%  '
$ws.Range("D4").Value = '''1.00'
$ws.Range("E4").Value = '  +0.02%  '
$ws.Range("D5").Value = '''357.76'
$ws.Range("E5").Value = '  +0.21%  '
$ws.Range("D6").Value = '''109.92'
$ws.Range("E6").Value = '  +0.68%  '
$ws.Range("D7").Value = '''0.558'
$ws.Range("E7").Value = '  +0.32%  '
$ws.Range("E8").Value = '  +0.06%  '
$ws.Range("D9").Value = '''0.635'
$ws.Range("E9").Value = '  +8.18%  '
$ws.Range("D10").Value = '''40.14'
$ws.Range("E10").Value = '  +0.86%  '
$ws.Range("E11").Value = '  +0.24%  '
$ws.Range("D12").Value = '''0.0841'
$ws.Range("E12").Value = '  -0.50%  '
$ws.Range("D13").Value = '''20.02'
$ws.Range("E13").Value = '  +2.74%  '
$ws.Range("D14").Value = '''7.82'
$ws.Range("E14").Value = '  +2.70%  '
$ws.Range("D15").Value = '3.254.42'
$ws.Range("E15").Value = '  +1.42%  '
$ws.Range("D16").Value = '2.841.27'
$ws.Range("E16").Value = '  +1.82%  '
$ws.Range("E17").Value = '  +1.07%  '
$ws.Range("D18").Value = '51.878.30'
$ws.Range("E18").Value = '  +0.45%  '
$ws.Range("D19").Value = '''7.70'
$ws.Range("E19").Value = '  +3.24%  '
$ws.Range("D20").Value = '''3.18'
$ws.Range("E20").Value = '  +3.55%  '
$ws.Range("D21").Value = '''13.70'
$ws.Range("E21").Value = '  +4.19%  '
$ws.Range("E22").Value = '  +1.03%  '
$ws.Range("D23").Value = '''70.51'
$ws.Range("E23").Value = '  +0.43%  '
$ws.Range("D24").Value = '''268.85'
$ws.Range("E24").Value = '  +0.03%  '
$ws.Range("E25").Value = '  +0.98%  '
$ws.Range("D26").Value = '''26.23'
$ws.Range("E26").Value = '  -0.62%  '
$ws.Range("E27").Value = '  +0.01%  '
$ws.Range("E28").Value = '  +0.58%  '
$ws.Range("D29").Value = '''10.40'
$ws.Range("E29").Value = '  +1.44%  '
$ws.Range("D30").Value = '''38.09'
$ws.Range("E30").Value = '  +9.41%  '
$ws.Range("B31").Value = 'Toncoin'
$ws.Range("C31").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D31").Value = '''2.24'
$ws.Range("E31").Value = '  +0.71%  '
$ws.Range("B32").Value = 'Filecoin'
$ws.Range("C32").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D32").Value = '''6.20'
$ws.Range("E32").Value = '  -0.35%  '
$ws.Range("B33").Value = 'OKB'
$ws.Range("C33").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D33").Value = '''52.15'
$ws.Range("E33").Value = '  +0.47%  '
$ws.Range("D34").Value = '''5.68'
$ws.Range("E34").Value = '  +10.59%  '
$ws.Range("D35").Value = '''0.0447'
$ws.Range("E35").Value = '  -1.22%  '
$ws.Range("D36").Value = '''0.0869'
$ws.Range("E36").Value = '  +3.50%  '
$ws.Range("D37").Value = '''1.00'
$ws.Range("E37").Value = '  +0.02%  '
$ws.Range("E38").Value = '  +1.10%  '
$ws.Range("E39").Value = '  +2.65%  '
$ws.Range("D40").Value = '''3.15'
$ws.Range("E40").Value = '  +0.40%  '
$ws.Range("E41").Value = '  +1.12%  '
$ws.Range("D42").Value = '''2.51'
$ws.Range("E42").Value = '  -1.09%  '
$ws.Range("D43").Value = '''22.03'
$ws.Range("E43").Value = '  +1.09%  '
$ws.Range("B44").Value = 'Monero'
$ws.Range("C44").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D44").Value = '''119.38'
$ws.Range("E44").Value = '  -0.19%  '
$ws.Range("B45").Value = 'WEMIXToken'
$ws.Range("C45").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D45").Value = '''2.20'
$ws.Range("E45").Value = '  -1.15%  '
$ws.Range("D46").Value = '''2.48'
$ws.Range("E46").Value = '  +8.48%  '
$ws.Range("B47").Value = 'Maker'
$ws.Range("C47").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D47").Value = '2.111.07'
$ws.Range("E47").Value = '  +1.21%  '
$ws.Range("B48").Value = 'NEARProtocol'
$ws.Range("C48").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D48").Value = '''3.40'
$ws.Range("E48").Value = '  +3.58%  '
$ws.Range("D49").Value = '''0.931'
$ws.Range("E49").Value = '  -1.05%  '
$ws.Range("E50").Value = '  +9.29%  '
$ws.Range("D51").Value = '''5.44'
$ws.Range("E51").Value = '  -4.50%  '
